# Auto-generated script to apply scheduled market-data updates to the Sheets workbook.
# For each changed cell, sets the new value directly; cells that no longer have data
# (removed entirely in the source data) are cleared via ClearContents().

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 216.73685
$ws.Range("I2").Value = 256
$ws.Range("J2").Value = 131.66667
$ws.Range("K2").Value = 256
$ws.Range("L2").Value = 131.66667
$ws.Range("M2").Value = -143
$ws.Range("N2").Value = -357.66667
$ws.Range("H9").Value = 8470.272000000001
$ws.Range("I9").Value = 9192.299999999999
$ws.Range("K9").Value = 9192.299999999999
$ws.Range("M9").Value = -9023.299999999999
$ws.Range("H28").Value = 1944.8
$ws.Range("I28").Value = 1922.8572
$ws.Range("J28").Value = 1996
$ws.Range("K28").Value = 1922.8572
$ws.Range("L28").Value = 1996
$ws.Range("M28").Value = -1437.8572
$ws.Range("N28").Value = -2966
$ws.Range("H41").Value = 1186.3572
$ws.Range("I41").Value = 1467.1818
$ws.Range("K41").Value = 1467.1818
$ws.Range("M41").Value = -1027.1818
$ws.Range("H86").Value = 3443.2173
$ws.Range("I86").Value = 2878.5
$ws.Range("J86").Value = 5476.2
$ws.Range("K86").Value = 2878.5
$ws.Range("L86").Value = 5476.2
$ws.Range("M86").Value = -1755.5
$ws.Range("N86").Value = -7722.2
$ws.Range("H89").Value = 3443.2173
$ws.Range("I89").Value = 2878.5
$ws.Range("J89").Value = 5476.2
$ws.Range("K89").Value = 14392.5
$ws.Range("L89").Value = 27381
$ws.Range("M89").Value = -8776.5
$ws.Range("N89").Value = -38613
$ws.Range("H92").Value = 6529.9375
$ws.Range("I92").Value = 7844.6924
$ws.Range("K92").Value = 7844.6924
$ws.Range("M92").Value = -6596.6924
$ws.Range("H99").Value = 5265.4287
$ws.Range("I99").Value = 339.5
$ws.Range("J99").Value = 11833.333
$ws.Range("K99").Value = 1018.5
$ws.Range("L99").Value = 35499.999
$ws.Range("M99").Value = 479.5
$ws.Range("N99").Value = -38495.999
$ws.Range("H132").Value = 10162749
$ws.Range("I132").Value = 11145631
$ws.Range("K132").Value = 33436893
$ws.Range("M132").Value = -33434363
$ws.Range("H137").Value = 1518.5883
$ws.Range("J137").Value = 2299.4
$ws.Range("L137").Value = 6898.200000000001
$ws.Range("N137").Value = -11998.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3275.5
$ws.Range("I32").Value = 2927.4614
$ws.Range("K32").Value = 2927.4614
$ws.Range("M32").Value = -2640.4614
$ws.Range("H45").Value = 20651.072
$ws.Range("I45").Value = 20950.125
$ws.Range("K45").Value = 20950.125
$ws.Range("M45").Value = -20573.125
$ws.Range("H74").Value = 53782.26
$ws.Range("I74").Value = 59874.47
$ws.Range("J74").Value = 1998.5
$ws.Range("K74").Value = 59874.47
$ws.Range("L74").Value = 1998.5
$ws.Range("M74").Value = -59000.47
$ws.Range("N74").Value = -3746.5
$ws.Range("H77").Value = 53782.26
$ws.Range("I77").Value = 59874.47
$ws.Range("J77").Value = 1998.5
$ws.Range("K77").Value = 299372.35
$ws.Range("L77").Value = 9992.5
$ws.Range("M77").Value = -295004.35
$ws.Range("N77").Value = -18728.5
$ws.Range("H97").Value = 1271.2
$ws.Range("I97").Value = 1599.6666
$ws.Range("J97").Value = 426.57144
$ws.Range("K97").Value = 1599.6666
$ws.Range("L97").Value = 426.57144
$ws.Range("M97").Value = -1103.6666
$ws.Range("N97").Value = -1418.57144
$ws.Range("H132").Value = 2139.0208
$ws.Range("I132").Value = 1933.9535
$ws.Range("K132").Value = 5801.860500000001
$ws.Range("M132").Value = -3271.860500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 20000
$ws.Range("J38").Value = 20000
$ws.Range("L38").Value = 20000
$ws.Range("N38").Value = -20832
$ws.Range("H86").Value = 13795.235
$ws.Range("I86").Value = 12822.786
$ws.Range("K86").Value = 12822.786
$ws.Range("M86").Value = -11699.786
$ws.Range("H89").Value = 13795.235
$ws.Range("I89").Value = 12822.786
$ws.Range("K89").Value = 64113.93
$ws.Range("M89").Value = -58497.93
$ws.Range("H110").Value = 36000
$ws.Range("J110").Value = 36000
$ws.Range("L110").Value = 36000
$ws.Range("N110").Value = -44180
$ws.Range("H132").Value = 72819.78
$ws.Range("J132").Value = 74322.375
$ws.Range("L132").Value = 74322.375
$ws.Range("N132").Value = -84442.375
$ws.Range("H133").Value = 65440
$ws.Range("J133").Value = 65440
$ws.Range("L133").Value = 65440
$ws.Range("N133").Value = -75560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 39469.895
$ws.Range("I31").Value = 45280.74
$ws.Range("J31").Value = 12740
$ws.Range("K31").Value = 45280.74
$ws.Range("L31").Value = 12740
$ws.Range("M31").Value = -44985.74
$ws.Range("N31").Value = -13330
$ws.Range("H34").Value = 39469.895
$ws.Range("I34").Value = 45280.74
$ws.Range("J34").Value = 12740
$ws.Range("K34").Value = 45280.74
$ws.Range("L34").Value = 12740
$ws.Range("M34").Value = -45078.74
$ws.Range("N34").Value = -13144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 799
$ws.Range("I35").Value = 594
$ws.Range("K35").Value = 1782
$ws.Range("M35").Value = -1494
$ws.Range("H47").Value = 100
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H87").Value = 8841.9
$ws.Range("I87").Value = 6534.875
$ws.Range("J87").Value = 18070
$ws.Range("K87").Value = 19604.625
$ws.Range("L87").Value = 54210
$ws.Range("M87").Value = -18356.625
$ws.Range("N87").Value = -56706
$ws.Range("H90").Value = 8841.9
$ws.Range("I90").Value = 6534.875
$ws.Range("J90").Value = 18070
$ws.Range("K90").Value = 58813.875
$ws.Range("L90").Value = 162630
$ws.Range("M90").Value = -52573.875
$ws.Range("N90").Value = -175110
$ws.Range("H92").Value = 713.3333
$ws.Range("I92").Value = 271.42856
$ws.Range("J92").Value = 1100
$ws.Range("K92").Value = 814.28568
$ws.Range("L92").Value = 3300
$ws.Range("M92").Value = 433.71432
$ws.Range("N92").Value = -5796
$ws.Range("H94").Value = 14999
$ws.Range("J94").Value = 14999
$ws.Range("L94").Value = 44997
$ws.Range("N94").Value = -46349
$ws.Range("H99").Value = 7782.5
$ws.Range("I99").Value = 2565
$ws.Range("K99").Value = 7695
$ws.Range("M99").Value = -5449
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H106").Value = 6962
$ws.Range("I106").Value = 6962
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 20886
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -19940
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1947.4445
$ws.Range("I97").Value = 1940.875
$ws.Range("K97").Value = 1940.875
$ws.Range("M97").Value = -1444.875
$ws.Range("H113").Value = 6403.8125
$ws.Range("I113").Value = 5252.2
$ws.Range("J113").Value = 7419.9414
$ws.Range("K113").Value = 5252.2
$ws.Range("L113").Value = 7419.9414
$ws.Range("M113").Value = -3082.2
$ws.Range("N113").Value = -11759.9414

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1958.6666
$ws.Range("I22").Value = 1789.5
$ws.Range("K22").Value = 1789.5
$ws.Range("M22").Value = -1494.5
$ws.Range("H23").Value = 514995
$ws.Range("I23").Value = 514995
$ws.Range("K23").Value = 514995
$ws.Range("M23").Value = -514765
$ws.Range("H27").Value = 1958.6666
$ws.Range("I27").Value = 1789.5
$ws.Range("K27").Value = 1789.5
$ws.Range("M27").Value = -1682.5
$ws.Range("H82").Value = 2800
$ws.Range("J82").Value = 3100
$ws.Range("L82").Value = 3100
$ws.Range("N82").Value = -3822
$ws.Range("H85").Value = 2800
$ws.Range("J85").Value = 3100
$ws.Range("L85").Value = 3100
$ws.Range("N85").Value = -5596
$ws.Range("H93").Value = 486163.4
$ws.Range("J93").Value = 4988.1113
$ws.Range("L93").Value = 4988.1113
$ws.Range("N93").Value = -7484.1113
$ws.Range("H132").Value = 3794.6428
$ws.Range("I132").Value = 3412.0417
$ws.Range("K132").Value = 10236.1251
$ws.Range("M132").Value = -7706.125100000001
$ws.Range("H141").Value = 84425
$ws.Range("J141").Value = 84425
$ws.Range("L141").Value = 84425
$ws.Range("N141").Value = -94785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 15001.857
$ws.Range("J45").Value = 15502.333
$ws.Range("L45").Value = 15502.333
$ws.Range("N45").Value = -16484.333
$ws.Range("H96").Value = 1665.3572
$ws.Range("I96").Value = 1208.75
$ws.Range("J96").Value = 1848
$ws.Range("K96").Value = 1208.75
$ws.Range("L96").Value = 1848
$ws.Range("M96").Value = 164.25
$ws.Range("N96").Value = -4594
